{"js": "// Locate the target paragraphs by scanning body.paragraphs:\n//  1) The second of four consecutive empty, strike-formatted paragraphs\n//     (between \"Launching of landing pods\" and \"Landing on new planet\")\n//     gets replaced with a brand-new, unformatted paragraph carrying the\n//     \"Monster sightings...\" text.\n//  2) The \"Centipede monster\" paragraph gains strikethrough formatting on\n//     both the paragraph mark and its run.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text,items/font/strikeThrough\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- Find paragraph 1: an empty, strikethrough paragraph that is the 2nd\n// of a run of 4 consecutive empty strikethrough paragraphs.\nlet targetIndex1 = -1;\nfor (let i = 0; i + 3 < items.length; i++) {\n  if (\n    items[i].text === \"\" && items[i].font.strikeThrough &&\n    items[i + 1].text === \"\" && items[i + 1].font.strikeThrough &&\n    items[i + 2].text === \"\" && items[i + 2].font.strikeThrough &&\n    items[i + 3].text === \"\" && items[i + 3].font.strikeThrough\n  ) {\n    targetIndex1 = i + 1; // the 2nd paragraph in the run of 4\n    break;\n  }\n}\nif (targetIndex1 === -1) {\n  throw new Error(\"Could not locate the 4-paragraph strikethrough run.\");\n}\n\n// --- Find paragraph 2: the paragraph whose text (ignoring the leading\n// tab) is \"Centipede monster\".\nlet targetIndex2 = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.replace(/^\\t+/, \"\") === \"Centipede monster\") {\n    targetIndex2 = i;\n    break;\n  }\n}\nif (targetIndex2 === -1) {\n  throw new Error('Could not locate the \"Centipede monster\" paragraph.');\n}\n\nconst pkgOpen =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>';\nconst pkgClose = \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\n// Edit 1: swap the empty strikethrough paragraph for a plain paragraph\n// with the new sentence (no paragraph- or run-level formatting at all).\nconst p1 = paragraphs.items[targetIndex1];\nconst range1 = p1.getRange();\nconst ooxml1 = pkgOpen +\n  \"<w:p><w:r><w:t>Monster sightings are everywhere. Documented and shared. Any persons going on expeditions is required to memorize.</w:t></w:r></w:p>\" +\n  pkgClose;\nrange1.insertOoxml(ooxml1, Word.InsertLocation.replace);\n\n// Edit 2: give the \"Centipede monster\" paragraph strikethrough formatting\n// on both the paragraph mark and the run.\nconst p2 = paragraphs.items[targetIndex2];\nconst range2 = p2.getRange();\nconst ooxml2 = pkgOpen +\n  \"<w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:tab/><w:t>Centipede monster</w:t></w:r></w:p>\" +\n  pkgClose;\nrange2.insertOoxml(ooxml2, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Locate the target paragraphs by scanning $d.Paragraphs:\n#  1) The second of four consecutive empty, strikethrough paragraphs\n#     (between \"Launching of landing pods\" and \"Landing on new planet\")\n#     gets replaced with a brand-new, unformatted paragraph carrying the\n#     \"Monster sightings...\" text.\n#  2) The \"Centipede monster\" paragraph gains strikethrough formatting on\n#     both the paragraph mark and its run.\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n# --- Find paragraph 1: an empty, strikethrough paragraph that is the 2nd\n# of a run of 4 consecutive empty strikethrough paragraphs.\n$targetIndex1 = -1\nfor ($i = 1; $i -le ($count - 3); $i++) {\n    $allMatch = $true\n    for ($j = 0; $j -le 3; $j++) {\n        $p = $d.Paragraphs.Item($i + $j)\n        $text = $p.Range.Text.Trim()\n        if (-not ($text -eq \"\" -and $p.Range.Font.StrikeThrough -eq -1)) {\n            $allMatch = $false\n        }\n    }\n    if ($allMatch) {\n        $targetIndex1 = $i + 1\n        break\n    }\n}\nif ($targetIndex1 -eq -1) {\n    throw \"Could not locate the 4-paragraph strikethrough run.\"\n}\n\n# --- Find paragraph 2: the paragraph whose trimmed text is\n# \"Centipede monster\" (ignoring the leading tab).\n$targetIndex2 = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.Trim()\n    if ($text -eq \"Centipede monster\") {\n        $targetIndex2 = $i\n        break\n    }\n}\nif ($targetIndex2 -eq -1) {\n    throw \"Could not locate the 'Centipede monster' paragraph.\"\n}\n\n$pkgOpen = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# Edit 1: swap the empty strikethrough paragraph for a plain paragraph\n# with the new sentence (no paragraph- or run-level formatting at all).\n$p1 = $d.Paragraphs.Item($targetIndex1)\n$xml1 = $pkgOpen + '<w:p><w:r><w:t>Monster sightings are everywhere. Documented and shared. Any persons going on expeditions is required to memorize.</w:t></w:r></w:p>' + $pkgClose\n$p1.Range.InsertXML($xml1)\n\n# Edit 2: give the \"Centipede monster\" paragraph strikethrough formatting\n# on both the paragraph mark and the run.\n$p2 = $d.Paragraphs.Item($targetIndex2)\n$xml2 = $pkgOpen + '<w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:tab/><w:t>Centipede monster</w:t></w:r></w:p>' + $pkgClose\n$p2.Range.InsertXML($xml2)\n"}
